# Auto-generated edit script applying the cryptos.xlsx diff (rows 2-51)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "'57.064.65"
$ws.Cells.Item(2, 5).Value = "  -7.50%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "'2.872.44"
$ws.Cells.Item(3, 5).Value = "  -5.63%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.03%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'549.45"
$ws.Cells.Item(5, 5).Value = "  -5.56%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'121.07"
$ws.Cells.Item(6, 5).Value = "  -7.23%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +0.16%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "'2.863.51"
$ws.Cells.Item(8, 5).Value = "  -5.91%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  -2.66%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  -10.40%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "'4.72"
$ws.Cells.Item(11, 5).Value = "  -10.52%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "'0.430"
$ws.Cells.Item(12, 5).Value = "  -2.25%  "

# Row 13
$ws.Cells.Item(13, 5).Value = "  -10.32%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "'31.28"
$ws.Cells.Item(14, 5).Value = "  -6.98%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "'0.118"
$ws.Cells.Item(15, 5).Value = "  -1.28%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "'3.345.53"
$ws.Cells.Item(16, 5).Value = "  -5.71%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "'2.870.20"
$ws.Cells.Item(17, 5).Value = "  -5.86%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "'57.033.25"
$ws.Cells.Item(18, 5).Value = "  -7.65%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "'6.41"
$ws.Cells.Item(19, 5).Value = "  +0.48%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "'406.63"
$ws.Cells.Item(20, 5).Value = "  -9.42%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "'12.73"
$ws.Cells.Item(21, 5).Value = "  -5.85%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "'0.649"
$ws.Cells.Item(22, 5).Value = "  -3.56%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  -8.51%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "'12.50"
$ws.Cells.Item(24, 5).Value = "  -2.83%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "'76.36"
$ws.Cells.Item(25, 5).Value = "  -5.73%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "'1.00"
$ws.Cells.Item(26, 5).Value = "  +0.03%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "'1.00"
$ws.Cells.Item(27, 5).Value = "  +0.01%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "'2.45"
$ws.Cells.Item(28, 5).Value = "  -4.38%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  -5.87%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "'7.06"
$ws.Cells.Item(30, 5).Value = "  -4.79%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "'6.00"
$ws.Cells.Item(31, 5).Value = "  -7.17%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "'24.49"
$ws.Cells.Item(32, 5).Value = "  -5.50%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "'0.0939"
$ws.Cells.Item(33, 5).Value = "  -3.58%  "

# Row 34
$ws.Cells.Item(34, 2).Value = "Filecoin"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(34, 4).Value = "'5.32"
$ws.Cells.Item(34, 5).Value = "  -7.21%  "

# Row 35
$ws.Cells.Item(35, 2).Value = "Stacks"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(35, 4).Value = "'2.00"
$ws.Cells.Item(35, 5).Value = "  -14.00%  "

# Row 36
$ws.Cells.Item(36, 2).Value = "OKB"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(36, 4).Value = "'48.27"
$ws.Cells.Item(36, 5).Value = "  -4.19%  "

# Row 37
$ws.Cells.Item(37, 2).Value = "Mantle"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(37, 4).Value = "'0.886"
$ws.Cells.Item(37, 5).Value = "  -9.09%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "'8.29"
$ws.Cells.Item(38, 5).Value = "  +4.63%  "

# Row 39
$ws.Cells.Item(39, 4).Value = "'0.0₃0612"
$ws.Cells.Item(39, 5).Value = "  -12.12%  "

# Row 40
$ws.Cells.Item(40, 4).Value = "'0.0341"
$ws.Cells.Item(40, 5).Value = "  -8.81%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  -3.70%  "

# Row 42
$ws.Cells.Item(42, 2).Value = "Bittensor"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(42, 4).Value = "'360.59"
$ws.Cells.Item(42, 5).Value = "  -5.33%  "

# Row 43
$ws.Cells.Item(43, 2).Value = "Maker"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(43, 4).Value = "'2.588.88"
$ws.Cells.Item(43, 5).Value = "  -4.05%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  +0.02%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "'2.32"
$ws.Cells.Item(45, 5).Value = "  -8.22%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "'117.65"
$ws.Cells.Item(46, 5).Value = "  -5.08%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  -5.43%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "'0.106"
$ws.Cells.Item(48, 5).Value = "  -2.02%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "'1.90"
$ws.Cells.Item(49, 5).Value = "  -5.45%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "'22.19"
$ws.Cells.Item(50, 5).Value = "  -7.27%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "'1.92"
$ws.Cells.Item(51, 5).Value = "  -6.88%  "

